# Update "想去人数" (want-to-go count) and "最低票价" (lowest price) figures
# across the four sheets to match the newly scraped data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "展览" (Exhibition)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 27073
$ws1.Range("F4").Value = 666
$ws1.Range("G4").Value = 60
$ws1.Range("F5").Value = 198
$ws1.Range("F6").Value = 578
$ws1.Range("F8").Value = 383
$ws1.Range("F9").Value = 492
$ws1.Range("F11").Value = 54
$ws1.Range("F12").Value = 319
$ws1.Range("F13").Value = 102
$ws1.Range("F14").Value = 518
$ws1.Range("F15").Value = 71
$ws1.Range("F16").Value = 1653
$ws1.Range("F17").Value = 271
$ws1.Range("F18").Value = 1033
$ws1.Range("F19").Value = 200
$ws1.Range("F23").Value = 122

# ---------------------------------------------------------------
# Sheet "演出" (Performance)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 4533
$ws2.Range("G3").Value = "不可售"
$ws2.Range("F8").Value = 7
$ws2.Range("F10").Value = 460
$ws2.Range("F16").Value = 79

# ---------------------------------------------------------------
# Sheet "本地生活" (Local Life)
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5204

# ---------------------------------------------------------------
# Sheet "全部类型" (All Types)
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5204
$ws4.Range("F5").Value = 27073
$ws4.Range("F6").Value = 4533
$ws4.Range("F7").Value = 666
$ws4.Range("G7").Value = 60
$ws4.Range("G8").Value = "不可售"
$ws4.Range("F10").Value = 198
$ws4.Range("F14").Value = 7
$ws4.Range("F16").Value = 460
$ws4.Range("F17").Value = 578
$ws4.Range("F21").Value = 383
$ws4.Range("F22").Value = 492
$ws4.Range("F24").Value = 54
$ws4.Range("F26").Value = 319
$ws4.Range("F27").Value = 102
$ws4.Range("F30").Value = 518
$ws4.Range("F31").Value = 71
$ws4.Range("F32").Value = 79
$ws4.Range("F33").Value = 1653
$ws4.Range("F34").Value = 271
$ws4.Range("F35").Value = 1034
$ws4.Range("F37").Value = 200
$ws4.Range("F42").Value = 122
